$wb = $excel.ActiveWorkbook

$wsData  = $wb.Worksheets.Item("Data")
$wsLogin = $wb.Worksheets.Item("Log-In creditials")

# --- "Data" sheet: was a 4-column placeholder row (swag/humble/boii/name),
#     becomes a 2-column header row for the receiver-email template ---
$wsData.Range("A1").Value = "Receiver Email"
$wsData.Range("B1").Value = "Subject"
$wsData.Range("C1:D1").ClearContents()

# Column A on "Data" gets an explicit best-fit-like width.
$wsData.Columns.Item(1).ColumnWidth = 12

# --- "Log-In creditials" sheet: fill in the actual credentials next to the
#     existing "Username:" / "Password:" labels (labels themselves unchanged) ---
$wsLogin.Range("B1").Value = "snirrfakturor"
$wsLogin.Range("B2").Value = "IrrSnirr96"

# --- Selection / active-sheet bookkeeping, matching the authored workbook ---
[void]$wsData.Range("C1:D1").Select()

[void]$wsLogin.Activate()
[void]$wsLogin.Range("A3").Select()
